$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep tracker_date/report_date columns as plain text (they look like dates
# "dd.mm.yyyy" and would otherwise get auto-converted to date serials).
$ws.Range("A156:B158").NumberFormat = "@"

# Row 156 - 11.03.2024 (report still dated 10.03.2024)
$ws.Range("A156").Value = "11.03.2024"
$ws.Range("B156").Value = "10.03.2024"
$ws.Range("C156").Value = 31045
$ws.Range("D156").Value = 12300
$ws.Range("E156").Value = 8400
$ws.Range("F156").Value = 72654
$ws.Range("G156").Value = 8663
$ws.Range("H156").Value = 6327
$ws.Range("I156").Value = 8000
$ws.Range("J156").Value = 425
$ws.Range("K156").Value = 113
$ws.Range("L156").Value = 4650
$ws.Range("M156").Value = "https://web.archive.org/web/20240311174900/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# Row 157 - 12.03.2024 (report still dated 10.03.2024)
$ws.Range("A157").Value = "12.03.2024"
$ws.Range("B157").Value = "10.03.2024"
$ws.Range("C157").Value = 31045
$ws.Range("D157").Value = 12300
$ws.Range("E157").Value = 8400
$ws.Range("F157").Value = 72654
$ws.Range("G157").Value = 8663
$ws.Range("H157").Value = 6327
$ws.Range("I157").Value = 8000
$ws.Range("J157").Value = 425
$ws.Range("K157").Value = 113
$ws.Range("L157").Value = 4650
$ws.Range("M157").Value = "https://web.archive.org/web/20240312034827/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# Row 158 - 13.03.2024 (report updated to 13.03.2024 with new figures)
$ws.Range("A158").Value = "13.03.2024"
$ws.Range("B158").Value = "13.03.2024"
$ws.Range("C158").Value = 31272
$ws.Range("D158").Value = 12300
$ws.Range("E158").Value = 8400
$ws.Range("F158").Value = 73024
$ws.Range("G158").Value = 8663
$ws.Range("H158").Value = 6327
$ws.Range("I158").Value = 8000
$ws.Range("J158").Value = 432
$ws.Range("K158").Value = 115
$ws.Range("L158").Value = 4650
$ws.Range("M158").Value = "https://web.archive.org/web/20240313155810/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# Restore default number format on the date-like text columns.
$ws.Range("A156:B158").NumberFormat = "General"
